$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Concatenate new data (Student ID "00075", Student Name "Carmen") onto the
# existing table, right after the last used row.
$ws.Range("A4").Value = "00075"
$ws.Range("B4").Value = "Carmen"

$ws.Range("B4").Select()
